$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CoursesPage")

# Fill in the missing test-case number for the row that was left blank
# (it sits between "3" and what used to be labelled "4"), and renumber
# every test case below it so the numbering runs continuously again
# (the old sheet restarted the count at 17 after the "Delete" section
# header instead of continuing on from 21).
$rowToNumber = @{
    5  = 4
    6  = 5
    7  = 6
    8  = 7
    9  = 8
    10 = 9
    11 = 10
    12 = 11
    13 = 12
    14 = 13
    16 = 14
    17 = 15
    18 = 16
    19 = 17
    20 = 18
    21 = 19
    22 = 20
    24 = 21
    25 = 22
    26 = 23
    27 = 24
    28 = 25
    29 = 26
}

foreach ($row in $rowToNumber.Keys) {
    $ws.Cells.Item($row, 1).Value2 = $rowToNumber[$row]
}

# Update the sheet's view/selection state to match where the author left
# the cursor after editing (selection on B35).
$ws.Range("B35").Select() | Out-Null
